$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.287.18"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "1.872.19"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "0.7079"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "241.63"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.07788"
$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("D9").Value = "0.3108"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").Value = "25.01"
$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("D11").Value = "0.08393"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "1.879.11"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "5.234"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").Value = "0.7167"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "6.118"
$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008332"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("D18").Value = "29.299.70"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").Value = "240.28"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("E23").Value = "  -2.43%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").Value = "4.403"
$ws.Range("E30").Value = "  -0.57%  "

$ws.Range("D31").Value = "4.313"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").Value = "1.254"
$ws.Range("E32").Value = "  -3.47%  "

$ws.Range("D33").Value = "0.05361"
$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("D34").Value = "1.936"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").Value = "1.239.61"
$ws.Range("E39").Value = "  +5.85%  "

$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").Value = "6.518"
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("D42").Value = "0.8902"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "2.023.44"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000126"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.06%  "

$ws.Range("D48").Value = "0.5198"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "1.788"
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").Value = "9.432"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").Value = "0.4332"
$ws.Range("E51").Value = "  +0.16%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.129.24"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.024"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.174"
$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7486"
$ws.Range("E36").Value = "  -3.65%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.26%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "72.33"
$ws.Range("E44").Value = "  -2.11%  "
